$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: new report week (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Crime Complaints table updates (rows 14-21 Murder..TOTAL, 22-30 Transit..Hate Crimes) ---
# Row 14
$ws.Range("G14").NumberFormat = "General"
$ws.Range("G14").Value = "'0"
$ws.Range("H14").NumberFormat = "General"
$ws.Range("H14").Value = "***.*"
# Row 15
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 2
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 300
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = -33.333333333333
# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 116.666666666667
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 25
$ws.Range("L16").Value = 138.095238095238
$ws.Range("M16").Value = 42.857142857142
$ws.Range("N16").Value = -66.216216216216
# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 138.461538461538
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 19
$ws.Range("K17").Value = 152.631578947368
$ws.Range("L17").Value = 182.352941176471
$ws.Range("M17").Value = 242.857142857143
$ws.Range("N17").Value = 29.729729729729
# Row 18
$ws.Range("C18").Value = 19
$ws.Range("E18").Value = 90
$ws.Range("F18").Value = 61
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = 48.780487804878
$ws.Range("I18").Value = 86
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = 50.877192982456
$ws.Range("L18").Value = 79.166666666666
$ws.Range("M18").Value = 36.507936507936
$ws.Range("N18").Value = -72.523961661341
# Row 19
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 31
$ws.Range("E19").Value = -25.806451612903
$ws.Range("F19").Value = 106
$ws.Range("G19").Value = 138
$ws.Range("H19").Value = -23.188405797101
$ws.Range("I19").Value = 170
$ws.Range("J19").Value = 215
$ws.Range("K19").Value = -20.930232558139
$ws.Range("L19").Value = 136.111111111111
$ws.Range("M19").Value = 161.538461538462
$ws.Range("N19").Value = 10.38961038961
# Row 20
$ws.Range("C20").Value = 10
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 36
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 63.636363636363
$ws.Range("I20").Value = 55
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = 66.666666666666
$ws.Range("L20").Value = 120
$ws.Range("M20").Value = 77.419354838709
$ws.Range("N20").Value = -89
# Row 21
$ws.Range("C21").Value = 69
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 23.214285714285
$ws.Range("F21").Value = 275
$ws.Range("G21").Value = 234
$ws.Range("H21").Value = 17.521367521367
$ws.Range("I21").Value = 413
$ws.Range("J21").Value = 353
$ws.Range("K21").Value = 16.99716713881
$ws.Range("L21").Value = 124.45652173913
$ws.Range("M21").Value = 94.811320754717
$ws.Range("N21").Value = -64.365832614322
# Row 22
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -50
# Row 23
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").NumberFormat = "General"
$ws.Range("F23").Value = "'0"
$ws.Range("H23").Value = -100
# Row 24
$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = 18.518518518518
$ws.Range("F24").Value = 233
$ws.Range("G24").Value = 212
$ws.Range("H24").Value = 9.905660377358
$ws.Range("I24").Value = 343
$ws.Range("J24").Value = 317
$ws.Range("K24").Value = 8.201892744479
$ws.Range("L24").Value = 89.502762430939
$ws.Range("M24").Value = 81.481481481481
# Row 25
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 64
$ws.Range("H25").Value = 42.222222222222
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 31.428571428571
$ws.Range("L25").Value = 135.897435897436
$ws.Range("M25").Value = 43.75
# Row 26
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 50
$ws.Range("L26").Value = 200
# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 3
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 11.111111111111
$ws.Range("L27").Value = 150
# Row 28
$ws.Range("F28").NumberFormat = "General"
$ws.Range("F28").Value = "'0"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").NumberFormat = "General"
$ws.Range("H28").Value = "***.*"
$ws.Range("N28").Value = -80
# Row 29
$ws.Range("F29").NumberFormat = "General"
$ws.Range("F29").Value = "'0"
$ws.Range("G29").NumberFormat = "General"
$ws.Range("G29").Value = "'0"
$ws.Range("H29").NumberFormat = "General"
$ws.Range("H29").Value = "***.*"
$ws.Range("N29").Value = -80
# Row 30
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H30").Value = -100
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("J30").Value = 1
$ws.Range("K30").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("K30").Value = -100
